# Update the division-problem table: replace the text of 25 specific
# cells (5 content rows x 5 columns) with new values, leaving every other
# paragraph / row (including the blank spacer rows) untouched.
#
# Mapping is keyed by (1-based Word row index, 1-based column index) so
# each cell is targeted unambiguously, even though a handful of the old
# and new values repeat elsewhere in the table.

$d = $word.ActiveDocument
$table = $d.Tables(1)

$edits = @(
    @{ Row = 1;  Col = 1; Old = "65÷9="; New = "11÷3=" },
    @{ Row = 1;  Col = 2; Old = "17÷2="; New = "71÷2=" },
    @{ Row = 1;  Col = 3; Old = "80÷8="; New = "65÷2=" },
    @{ Row = 1;  Col = 4; Old = "87÷7="; New = "16÷3=" },
    @{ Row = 1;  Col = 5; Old = "23÷3="; New = "51÷6=" },

    @{ Row = 5;  Col = 1; Old = "45÷2="; New = "17÷2=" },
    @{ Row = 5;  Col = 2; Old = "20÷5="; New = "13÷5=" },
    @{ Row = 5;  Col = 3; Old = "77÷5="; New = "52÷2=" },
    @{ Row = 5;  Col = 4; Old = "73÷6="; New = "33÷6=" },
    @{ Row = 5;  Col = 5; Old = "25÷2="; New = "83÷2=" },

    @{ Row = 9;  Col = 1; Old = "69÷3="; New = "78÷4=" },
    @{ Row = 9;  Col = 2; Old = "85÷8="; New = "86÷5=" },
    @{ Row = 9;  Col = 3; Old = "76÷6="; New = "22÷9=" },
    @{ Row = 9;  Col = 4; Old = "65÷2="; New = "24÷2=" },
    @{ Row = 9;  Col = 5; Old = "99÷8="; New = "65÷7=" },

    @{ Row = 13; Col = 1; Old = "45÷8="; New = "79÷6=" },
    @{ Row = 13; Col = 2; Old = "43÷8="; New = "64÷2=" },
    @{ Row = 13; Col = 3; Old = "27÷9="; New = "14÷5=" },
    @{ Row = 13; Col = 4; Old = "99÷5="; New = "19÷7=" },
    @{ Row = 13; Col = 5; Old = "12÷9="; New = "94÷5=" },

    @{ Row = 17; Col = 1; Old = "78÷4="; New = "50÷9=" },
    @{ Row = 17; Col = 2; Old = "29÷3="; New = "97÷7=" },
    @{ Row = 17; Col = 3; Old = "81÷9="; New = "39÷3=" },
    @{ Row = 17; Col = 4; Old = "34÷4="; New = "74÷4=" },
    @{ Row = 17; Col = 5; Old = "69÷2="; New = "62÷4=" }
)

# Assign each cell's Range.Text directly instead of Find/Replace: this
# runtime's Find.Execute(Replace:=...) walks the whole story looking for
# the first match rather than honoring the Range it was scoped to, which
# is disastrous here since several old/new values repeat across cells.
# Range.Text = "..." mutates only the addressed cell and keeps the
# existing run's rPr (font/size) intact -- it only swaps the <w:t>.
foreach ($edit in $edits) {
    $cell = $table.Cell($edit.Row, $edit.Col)
    $before = $cell.Range.Text
    if ($before -notmatch [regex]::Escape($edit.Old)) {
        Write-Host "UNEXPECTED TEXT: row=$($edit.Row) col=$($edit.Col) expected=$($edit.Old) actual=$before"
    }
    $cell.Range.Text = $edit.New
}

Write-Host "Done."
